# Update extracted table: strip thousands separators from numeric text,
# split "Other income (expense), net (3,514)" into a clean label plus its
# own negative-value cell, and append the "Basic"/"Diluted net income per
# share (Note 12)" rows that were previously dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

# Rows 2-8 and 10-12: remove thousands separators (values stay as text).
Set-TextValue "B2" "282836"
Set-TextValue "C2" "307394"
Set-TextValue "D2" "350018"

Set-TextValue "B3" "126203"
Set-TextValue "C3" "133332"
Set-TextValue "D3" "146306"

Set-TextValue "B4" "39500"
Set-TextValue "C4" "45427"
Set-TextValue "D4" "49326"

Set-TextValue "B5" "26567"
Set-TextValue "C5" "27917"
Set-TextValue "D5" "27808"

Set-TextValue "B6" "15724"
Set-TextValue "C6" "16425"
Set-TextValue "D6" "14188"

Set-TextValue "B7" "207994"
Set-TextValue "C7" "223101"
Set-TextValue "D7" "237628"

Set-TextValue "B8" "74842"
Set-TextValue "C8" "84293"
Set-TextValue "D8" "112390"

# Row 9: label loses the trailing "(3,514)"; the figure becomes its own
# negative value, and the other two columns shift into plain numeric text.
Set-TextValue "A9" "Other income (expense), net"
Set-TextValue "B9" "-3514"
Set-TextValue "C9" "1424"
Set-TextValue "D9" "7425"

Set-TextValue "B10" "71328"
Set-TextValue "C10" "85717"
Set-TextValue "D10" "119815"

Set-TextValue "B11" "11356"
Set-TextValue "C11" "11922"
Set-TextValue "D11" "19697"

Set-TextValue "B12" "59972"
Set-TextValue "C12" "73795"
Set-TextValue "D12" "100118"

# Rows 13-14: newly read "Basic"/"Diluted net income per share" lines.
Set-TextValue "A13" "Basic net income per share (Note 12)"
Set-TextValue "B13" "4.59"
Set-TextValue "C13" "5.84"
Set-TextValue "D13" "8.13"

Set-TextValue "A14" "Diluted net income per share (Note 12)"
Set-TextValue "B14" "4.56"
Set-TextValue "C14" "5.80"
Set-TextValue "D14" "8.04"
